$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 163; existing rows 163-196 shift down to 164-197.
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with the new weekly price record.
$ws.Cells.Item(163, 1).Value = 3
$ws.Cells.Item(163, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(163, 3).Value = "Coquimbo"
$ws.Cells.Item(163, 4).Value = 44476
$ws.Cells.Item(163, 5).Value = 5
$ws.Cells.Item(163, 6).Value = 100112012
$ws.Cells.Item(163, 7).Value = "Espinaca"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 220
$ws.Cells.Item(163, 11).Value = 2500
$ws.Cells.Item(163, 12).Value = 2500
$ws.Cells.Item(163, 13).Value = 2500
$ws.Cells.Item(163, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(163, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(163, 16).Value = 833
$ws.Cells.Item(163, 17).Value = 3
$ws.Cells.Item(163, 18).Value = "Hortaliza"
